$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = "bcspwr05.mtx"
$ws.Cells.Item(2, 3).Value = "MP"
$ws.Cells.Item(2, 4).Value = 5.177624323644659
$ws.Cells.Item(2, 5).Value = 81
$ws.Cells.Item(2, 6).Value = 0.009785652160644531
$ws.Cells.Item(2, 7).Value = 443
$ws.Cells.Item(2, 8).Value = "pattern"
$ws.Cells.Item(2, 9).Value = "symmetric"

# Row 3
$ws.Cells.Item(3, 2).Value = "bcspwr05.mtx"
$ws.Cells.Item(3, 3).Value = "MP_Aitken"
$ws.Cells.Item(3, 4).Value = 5.179112393033175
$ws.Cells.Item(3, 5).Value = 34
$ws.Cells.Item(3, 6).Value = 0.001904010772705078
$ws.Cells.Item(3, 7).Value = 443
$ws.Cells.Item(3, 8).Value = "pattern"
$ws.Cells.Item(3, 9).Value = "symmetric"

# Row 4
$ws.Cells.Item(4, 2).Value = "bcsstk01.mtx"
$ws.Cells.Item(4, 3).Value = "MP"
$ws.Cells.Item(4, 4).Value = 3014181127.26153
$ws.Cells.Item(4, 5).Value = 134
$ws.Cells.Item(4, 6).Value = 0.003325700759887695
$ws.Cells.Item(4, 7).Value = 48
$ws.Cells.Item(4, 8).Value = "real"
$ws.Cells.Item(4, 9).Value = "symmetric"

# Row 5
$ws.Cells.Item(5, 2).Value = "bcsstk01.mtx"
$ws.Cells.Item(5, 3).Value = "MP_Aitken"
$ws.Cells.Item(5, 4).Value = 3015647293.232991
$ws.Cells.Item(5, 5).Value = 85
$ws.Cells.Item(5, 6).Value = 0.002544403076171875
$ws.Cells.Item(5, 7).Value = 48
$ws.Cells.Item(5, 8).Value = "real"
$ws.Cells.Item(5, 9).Value = "symmetric"

# Row 6
$ws.Cells.Item(6, 2).Value = "bcspwr10.mtx"
$ws.Cells.Item(6, 3).Value = "MP"
$ws.Cells.Item(6, 4).Value = 6.809359881800338
$ws.Cells.Item(6, 5).Value = 186
$ws.Cells.Item(6, 6).Value = 2.760515213012695
$ws.Cells.Item(6, 7).Value = 5300
$ws.Cells.Item(6, 8).Value = "pattern"
$ws.Cells.Item(6, 9).Value = "symmetric"

# Row 7
$ws.Cells.Item(7, 2).Value = "bcspwr10.mtx"
$ws.Cells.Item(7, 3).Value = "MP_Aitken"
$ws.Cells.Item(7, 4).Value = 6.836792424603274
$ws.Cells.Item(7, 5).Value = 81
$ws.Cells.Item(7, 6).Value = 1.163532972335815
$ws.Cells.Item(7, 7).Value = 5300
$ws.Cells.Item(7, 8).Value = "pattern"
$ws.Cells.Item(7, 9).Value = "symmetric"

# Row 8
$ws.Cells.Item(8, 2).Value = "bcspwr06.mtx"
$ws.Cells.Item(8, 3).Value = "MP"
$ws.Cells.Item(8, 4).Value = 5.618009280534263
$ws.Cells.Item(8, 5).Value = 124
$ws.Cells.Item(8, 6).Value = 0.1332681179046631
$ws.Cells.Item(8, 7).Value = 1454
$ws.Cells.Item(8, 8).Value = "pattern"
$ws.Cells.Item(8, 9).Value = "symmetric"

# Row 9
$ws.Cells.Item(9, 2).Value = "bcspwr06.mtx"
$ws.Cells.Item(9, 3).Value = "MP_Aitken"
$ws.Cells.Item(9, 4).Value = 5.622463254476138
$ws.Cells.Item(9, 5).Value = 49
$ws.Cells.Item(9, 6).Value = 0.05295443534851074
$ws.Cells.Item(9, 7).Value = 1454
$ws.Cells.Item(9, 8).Value = "pattern"
$ws.Cells.Item(9, 9).Value = "symmetric"

# Row 10
$ws.Cells.Item(10, 2).Value = "bcspwr07.mtx"
$ws.Cells.Item(10, 3).Value = "MP"
$ws.Cells.Item(10, 4).Value = 5.664183085589989
$ws.Cells.Item(10, 5).Value = 99
$ws.Cells.Item(10, 6).Value = 0.131089448928833
$ws.Cells.Item(10, 7).Value = 1612
$ws.Cells.Item(10, 8).Value = "pattern"
$ws.Cells.Item(10, 9).Value = "symmetric"

# Row 11
$ws.Cells.Item(11, 2).Value = "bcspwr07.mtx"
$ws.Cells.Item(11, 3).Value = "MP_Aitken"
$ws.Cells.Item(11, 4).Value = 5.667241609924945
$ws.Cells.Item(11, 5).Value = 42
$ws.Cells.Item(11, 6).Value = 0.05905938148498535
$ws.Cells.Item(11, 7).Value = 1612
$ws.Cells.Item(11, 8).Value = "pattern"
$ws.Cells.Item(11, 9).Value = "symmetric"

# Row 12
$ws.Cells.Item(12, 2).Value = "bcsstk03.mtx"
$ws.Cells.Item(12, 3).Value = "MP"
$ws.Cells.Item(12, 4).Value = 199732738576.2932
$ws.Cells.Item(12, 5).Value = 25
$ws.Cells.Item(12, 6).Value = 0.001411199569702148
$ws.Cells.Item(12, 7).Value = 112
$ws.Cells.Item(12, 8).Value = "real"
$ws.Cells.Item(12, 9).Value = "symmetric"

# Row 13
$ws.Cells.Item(13, 2).Value = "bcsstk03.mtx"
$ws.Cells.Item(13, 3).Value = "MP_Aitken"
$ws.Cells.Item(13, 4).Value = 199734791972.886
$ws.Cells.Item(13, 5).Value = 20
$ws.Cells.Item(13, 6).Value = 0.001338720321655273
$ws.Cells.Item(13, 7).Value = 112
$ws.Cells.Item(13, 8).Value = "real"
$ws.Cells.Item(13, 9).Value = "symmetric"

# Row 14
$ws.Cells.Item(14, 2).Value = "bcsstk02.mtx"
$ws.Cells.Item(14, 3).Value = "MP"
$ws.Cells.Item(14, 4).Value = 18225.28861825787
$ws.Cells.Item(14, 5).Value = 28
$ws.Cells.Item(14, 6).Value = 0.001276731491088867
$ws.Cells.Item(14, 7).Value = 66
$ws.Cells.Item(14, 8).Value = "real"
$ws.Cells.Item(14, 9).Value = "symmetric"

# Row 15
$ws.Cells.Item(15, 2).Value = "bcsstk02.mtx"
$ws.Cells.Item(15, 3).Value = "MP_Aitken"
$ws.Cells.Item(15, 4).Value = 18224.71635802562
$ws.Cells.Item(15, 5).Value = 17
$ws.Cells.Item(15, 6).Value = 0.000827789306640625
$ws.Cells.Item(15, 7).Value = 66
$ws.Cells.Item(15, 8).Value = "real"
$ws.Cells.Item(15, 9).Value = "symmetric"

# Row 16
$ws.Cells.Item(16, 2).Value = "bcsstk05.mtx"
$ws.Cells.Item(16, 3).Value = "MP"
$ws.Cells.Item(16, 4).Value = 6197043.861659037
$ws.Cells.Item(16, 5).Value = 38
$ws.Cells.Item(16, 6).Value = 0.002539157867431641
$ws.Cells.Item(16, 7).Value = 153
$ws.Cells.Item(16, 8).Value = "real"
$ws.Cells.Item(16, 9).Value = "symmetric"

# Row 17
$ws.Cells.Item(17, 2).Value = "bcsstk05.mtx"
$ws.Cells.Item(17, 3).Value = "MP_Aitken"
$ws.Cells.Item(17, 4).Value = 6196986.835041617
$ws.Cells.Item(17, 5).Value = 27
$ws.Cells.Item(17, 6).Value = 0.001700401306152344
$ws.Cells.Item(17, 7).Value = 153
$ws.Cells.Item(17, 8).Value = "real"
$ws.Cells.Item(17, 9).Value = "symmetric"

# Row 18
$ws.Cells.Item(18, 2).Value = "bcspwr08.mtx"
$ws.Cells.Item(18, 3).Value = "MP"
$ws.Cells.Item(18, 4).Value = 5.783741344069723
$ws.Cells.Item(18, 5).Value = 152
$ws.Cells.Item(18, 6).Value = 0.209707498550415
$ws.Cells.Item(18, 7).Value = 1624
$ws.Cells.Item(18, 8).Value = "pattern"
$ws.Cells.Item(18, 9).Value = "symmetric"

# Row 19
$ws.Cells.Item(19, 2).Value = "bcspwr08.mtx"
$ws.Cells.Item(19, 3).Value = "MP_Aitken"
$ws.Cells.Item(19, 4).Value = 5.794950776491883
$ws.Cells.Item(19, 5).Value = 57
$ws.Cells.Item(19, 6).Value = 0.07850742340087891
$ws.Cells.Item(19, 7).Value = 1624
$ws.Cells.Item(19, 8).Value = "pattern"
$ws.Cells.Item(19, 9).Value = "symmetric"

# Row 20
$ws.Cells.Item(20, 2).Value = "bcsstk04.mtx"
$ws.Cells.Item(20, 3).Value = "MP"
$ws.Cells.Item(20, 4).Value = 9549224.597917093
$ws.Cells.Item(20, 5).Value = 18
$ws.Cells.Item(20, 6).Value = 0.001143455505371094
$ws.Cells.Item(20, 7).Value = 132
$ws.Cells.Item(20, 8).Value = "real"
$ws.Cells.Item(20, 9).Value = "symmetric"

# Row 21
$ws.Cells.Item(21, 2).Value = "bcsstk04.mtx"
$ws.Cells.Item(21, 3).Value = "MP_Aitken"
$ws.Cells.Item(21, 4).Value = 9549371.073998552
$ws.Cells.Item(21, 5).Value = 14
$ws.Cells.Item(21, 6).Value = 0.0007519721984863281
$ws.Cells.Item(21, 7).Value = 132
$ws.Cells.Item(21, 8).Value = "real"
$ws.Cells.Item(21, 9).Value = "symmetric"
